# Adds a new weekly price record for "Berenjena" (Vega Modelo de Temuco)
# at row 166, pushing the existing rows 166:238 down to 167:239 (and
# extending the sheet's used range from R238 to R239 accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 166; Excel shifts rows
# 166..238 down to 167..239 and inherits formatting (e.g. the date
# number format on column D) from the row being pushed down.
$ws.Rows(166).Insert()

# Populate the newly inserted row 166 with the new record. Most fields
# mirror the record that used to sit at row 166 (now at row 167) - only
# the date (D) and volume (J) are new values.
$ws.Cells.Item(166, 1).Value = 10
$ws.Cells.Item(166, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(166, 3).Value = "La Araucanía"
$ws.Cells.Item(166, 4).Value = 44609
$ws.Cells.Item(166, 5).Value = 9
$ws.Cells.Item(166, 6).Value = 100112001
$ws.Cells.Item(166, 7).Value = "Berenjena"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 150
$ws.Cells.Item(166, 11).Value = 12000
$ws.Cells.Item(166, 12).Value = 12000
$ws.Cells.Item(166, 13).Value = 12000
$ws.Cells.Item(166, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(166, 15).Value = "Región del Maule"
$ws.Cells.Item(166, 16).Value = 200
$ws.Cells.Item(166, 17).Value = 60
$ws.Cells.Item(166, 18).Value = "Hortaliza"
